# "icons on each button vers2 incl mvn clean package update"
#
# Adds a new "icons" worksheet at the end of the workbook (after "comments"),
# makes it the active sheet, and populates a small block of cells (A2:D6)
# mirroring the layout already used by the "comments" sheet - two labelled
# cells ("hdd.png") plus a handful of formatted-but-empty neighbour cells.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end of
# the tab strip (Excel's default Add() would insert before the active
# sheet, so we explicitly anchor it After: the last worksheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$iconsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$iconsSheet.Name = "icons"

# New shared-string value used twice in the header row.
$iconsSheet.Range("A2").Value = "hdd.png"
$iconsSheet.Range("B2").Value = "hdd.png"

# Touch the remaining cells that carry formatting in the target layout so
# they materialize in the sheet (value stays blank, matching the source
# "comments" sheet pattern).
$iconsSheet.Range("A2:D3").Font.Size = 10
$iconsSheet.Range("B4").Font.Size = 10
$iconsSheet.Range("D4").Font.Size = 10
$iconsSheet.Range("B5").Font.Size = 10
$iconsSheet.Range("D5").Font.Size = 10
$iconsSheet.Range("D6").Font.Size = 10

Write-Output "Added 'icons' sheet with $($wb.Worksheets.Count) total sheets; active tab is now $($iconsSheet.Name)"
